$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 245 (old row 245 shifts down to 246)
$ws.Rows.Item(245).Insert()

# Insert a new row at 247 (after the shifted row 246) for the new trailing record
$ws.Rows.Item(247).Insert()

# New row 245: 2024-06-11 data
$ws.Cells.Item(245, 1).Value = 45454.2916666667
$ws.Cells.Item(245, 2).Value = 803152
$ws.Cells.Item(245, 3).Value = 3.17499995231628
$ws.Cells.Item(245, 4).Value = 3.07999992370605
$ws.Cells.Item(245, 5).Value = 3.13499999046326
$ws.Cells.Item(245, 6).Value = 3.07999992370605
$g245 = $ws.Cells.Item(245, 7)
$g245.NumberFormat = "@"
$g245.Value = "3.07999992370605"
$g245.Style = "Normal"
$ws.Cells.Item(245, 8).Value = "YACHT.MI"

# Row 246 (previously row 245): only the date value changes, other columns unchanged
$ws.Cells.Item(246, 1).Value = 45455.2916666667

# New row 247: 2024-06-13 data
$ws.Cells.Item(247, 1).Value = 45456.6494560185
$ws.Cells.Item(247, 2).Value = 326171
$ws.Cells.Item(247, 3).Value = 3.10999989509583
$ws.Cells.Item(247, 4).Value = 3.01999998092651
$ws.Cells.Item(247, 5).Value = 3.10999989509583
$ws.Cells.Item(247, 6).Value = 3.01999998092651
$g247 = $ws.Cells.Item(247, 7)
$g247.NumberFormat = "@"
$g247.Value = "3.01999998092651"
$g247.Style = "Normal"
$ws.Cells.Item(247, 8).Value = "YACHT.MI"
